$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert Opt Portfolio (C) and Opt Portfolio with View (D) columns
# back to the previous working commit's values.

$ws.Range("C2").Value = 0.2454406141024289
$ws.Range("D2").Value = 0.245440628078421

$ws.Range("C3").Value = 0.1815773618901335
$ws.Range("D3").Value = 0.1815773450933061

$ws.Range("C4").Value = 0.09883777827451906
$ws.Range("D4").Value = 0.09883777827451899

$ws.Range("C5").Value = 0.09883777827451906
$ws.Range("D5").Value = 0.09883777827451901

$ws.Range("C6").Value = 0.09883777827451906
$ws.Range("D6").Value = 0.09883777827451899

$ws.Range("C7").Value = 0.11995695347692
$ws.Range("D7").Value = 0.1199569551788409

$ws.Range("C8").Value = 0.1565117357069661
$ws.Range("D8").Value = 0.1565117368258783
